$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "grauAcima" table (rows 107-120) had duplicated content from the
# "ordenado" table by mistake. Fix the Nº inputs / Tipo de entradas / Valor
# especifico / Casos de Teste cells so they refer to "ficheiros" instead of
# "ordenado", matching the other maximoFicheiros-style tables.

# Row 109: Nº inputs / Tipo de entradas
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = "≠ 1"

# Row 110-111: Valor especifico (Valida / Invalida)
# NOTE: set B111 before B110 so any newly introduced shared strings land in
# the same order as the reference workbook.
$ws.Range("B111").Value = "numFiles > 0, ficheiros ≠ null"
$ws.Range("B110").Value = "numFiles = int"
$ws.Range("C110").Value = "numFiles ≠ int"
$ws.Range("C111").Value = "numFiles =< 0"

# Rows 114-120: Casos de Teste
$ws.Range("B114").Value = "numFiles = 1"
$ws.Range("B115").Value = "numFiles = 1 2"
$ws.Range("B116").Value = "numFiles ="
$ws.Range("B117").Value = "numFiles = ""a"""
$ws.Range("B118").Value = "numFiles = 1"
$ws.Range("B119").Value = "numFiles = 0"
$ws.Range("B120").Value = "numFiles = 1"
